$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting so date-like strings are not auto-converted to dates.
$ws.Range("D17:D20").NumberFormat = "@"

# New row inserted at row 17 (existing rows 17-19 shift down to 18-20).
$ws.Range("A17").Value = "U7sA0AEACAAJ"
$ws.Range("B17").Value = "El Señor de los Anillos III"
$ws.Range("C17").Value = "Los ejércitos del Señor Oscuro van extendiendo cada vez más su maléfica sombra por la Tierra Media. Hombres, elfos y enanos unen sus fuerzas para presentar batalla a Sauron y sus huestes. Ajenos a estos preparativos, Frodo y Sam siguen adentrándose en el país de Mordor en su heroico viaje para destruir el Anillo de Poder en las Grietas del Destino."
$ws.Range("D17").Value = "2010-04-28"
$ws.Range("E17").Value = "John Ronald Reuel Tolkien"

# Former row 17 data moves to row 18.
$ws.Range("A18").Value = "UfYGAAAACAAJ"
$ws.Range("B18").Value = "El señor de los anillos"
$ws.Range("C18").Value = "Desconocido"
$ws.Range("D18").Value = "2002-02"
$ws.Range("E18").Value = "Kurt D. Bruner, Jim Ware"

# Former row 18 data moves to row 19.
$ws.Range("A19").Value = "WmdWtQAACAAJ"
$ws.Range("B19").Value = "El Señor de los anillos"
$ws.Range("C19").Value = "Desconocido"
$ws.Range("D19").Value = "2002"
$ws.Range("E19").Value = "J. R. R. Tolkien"

# Former row 19 data moves to row 20.
$ws.Range("A20").Value = "ZVwX0QEACAAJ"
$ws.Range("B20").Value = "El Señor de los Anillos"
$ws.Range("C20").Value = "Desconocido"
$ws.Range("D20").Value = "1985"
$ws.Range("E20").Value = "J. R. R. Tolkien"
